$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '98.462.43'
$ws.Range('E2').Value = '  -0.58%  '
$ws.Range('D3').Value = '3.387.64'
$ws.Range('E3').Value = '  +0.29%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '258.72'
$ws.Range('D5').Style = 'Normal'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '668.86'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +6.33%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.56'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +12.89%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.458'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +16.58%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '1.10'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +27.67%  '
$ws.Range('E10').Value = '  +0.00%  '
$ws.Range('D11').Value = '3.388.11'
$ws.Range('E11').Value = '  +0.33%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.210'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +5.42%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '42.55'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +17.78%  '
$ws.Range('E14').Value = '  +7.75%  '
$ws.Range('D15').Value = '98.424.53'
$ws.Range('E15').Value = '  -0.38%  '
$ws.Range('D16').Value = '4.022.43'
$ws.Range('E16').Value = '  +0.49%  '
$ws.Range('E17').Value = '  +2.41%  '
$ws.Range('D18').Value = '3.381.14'
$ws.Range('E18').Value = '  -0.04%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.63'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +24.78%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '17.03'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +11.42%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '3.59'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.02%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '530.08'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +7.42%  '
$ws.Range('E23').Value = '  +12.79%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.0000216'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.57%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.441'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +57.17%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '6.31'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +11.89%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '101.71'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +14.86%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '12.74'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +6.43%  '
$ws.Range('D29').Value = '3.559.22'
$ws.Range('E29').Value = '  +0.04%  '
$ws.Range('E30').Value = '  +15.50%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.999'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.05%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '11.19'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +16.36%  '
$ws.Range('E33').Value = '  -0.56%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.00'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.52%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '29.83'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +6.18%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.545'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +18.22%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '7.94'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +8.45%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.14'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +8.95%  '
$ws.Range('E39').Value = '  +7.32%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '527.48'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +5.39%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0454'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +38.98%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.34'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +5.71%  '
$ws.Range('E43').Value = '  -0.85%  '
$ws.Range('E44').Value = '  -1.33%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.840'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +6.83%  '
$ws.Range('E46').Value = '  +2.56%  '
$ws.Range('E47').Value = '  +0.03%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.08'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +6.95%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.86'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +19.18%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.53'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +12.25%  '
$ws.Range('B51').Value = 'Filecoin'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '5.12'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +10.76%  '
